$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update loading_percent values for data rows (sheet rows 2-25) per the recorded diff.
# Columns changed: B, C, E, F, G, H, K, O (columns D, I, J, L, M, N remain 0 / unchanged).

$colB = @(8.927621164094839, 8.526629530229632, 8.270977669816366, 8.164556964605778, 8.146754852750322, 8.26955132392925, 8.791393141476416, 9.734860899471803, 10.37384732087933, 10.65189366260272, 10.75530804548531, 10.73312017181304, 10.66043952336064, 10.61567466038114, 10.35541685676092, 10.19247639909303, 10.09757124967997, 10.06523625394151, 10.20994490776526, 10.6818389223658, 10.97929488734697, 10.82155644933733, 10.20205121692489, 9.488801940317916)
for ($i = 0; $i -lt $colB.Count; $i++) {
    $ws.Cells.Item(2 + $i, 2).Value = $colB[$i]
}

$colC = @(7.036167523550205, 6.910838680868044, 6.832580812537681, 6.8003922935992, 6.795030337938435, 6.832147871703977, 6.993243038473216, 7.297509743673445, 7.512327934959393, 7.607840768391873, 7.643666600008879, 7.635966497066062, 7.610795190189047, 7.595331673643363, 7.506039311388722, 7.450677186303041, 7.418627841179493, 7.407741805082265, 7.456592139469163, 7.618198109836638, 7.721808212091203, 7.666701272155782, 7.453918675375439, 7.216617112978913)
for ($i = 0; $i -lt $colC.Count; $i++) {
    $ws.Cells.Item(2 + $i, 3).Value = $colC[$i]
}

$colE = @(16.52273998434378, 15.58558963047363, 14.98501082362067, 14.73421049413865, 14.69220819127422, 14.98165258440299, 16.20497880181815, 18.48836619208137, 20.12722409329579, 20.83039812197094, 21.09061763039483, 21.03484358170894, 20.85192774689739, 20.7390986531409, 20.08042061732207, 19.66551932813181, 19.42289089463626, 19.34005509982445, 19.71009880299404, 20.90581869343055, 21.65199958369428, 21.25696719023441, 19.68995718400925, 17.84706862855227)
for ($i = 0; $i -lt $colE.Count; $i++) {
    $ws.Cells.Item(2 + $i, 5).Value = $colE[$i]
}

$colF = @(16.86991607391245, 15.89584955866815, 15.26997757108489, 15.008197319934, 14.96433081551589, 15.26647399323133, 16.5399640634477, 19.00274580682531, 20.67494806633232, 21.3917225636224, 21.65686569030329, 21.60004134736742, 21.4136618050453, 21.29868154950795, 20.62722412089977, 20.20408069597325, 19.95656407809801, 19.87204792380568, 20.24955283636154, 21.46857628470577, 22.22866616901552, 21.82633154458858, 20.22900810905287, 18.34778573295695)
for ($i = 0; $i -lt $colF.Count; $i++) {
    $ws.Cells.Item(2 + $i, 6).Value = $colF[$i]
}

$colG = @(21.38898576005713, 21.57569643755049, 21.70021256605301, 21.75341715725123, 21.76239993803712, 21.70092015427344, 21.45130199178301, 21.04105795843175, 20.78922257744759, 20.68569513983817, 20.64810182778598, 20.65612622360602, 20.68256990676834, 20.69897783800486, 20.79621238184823, 20.85870537834227, 20.89568681476896, 20.90838550318996, 20.85194538379643, 20.67475886389908, 20.56835430316271, 20.62427671992124, 20.85499829781725, 21.14342660927449)
for ($i = 0; $i -lt $colG.Count; $i++) {
    $ws.Cells.Item(2 + $i, 7).Value = $colG[$i]
}

$colH = @(12.58980090402059, 12.65492127157398, 12.69722219753125, 12.71504323456757, 12.71803763513342, 12.69746017679272, 12.61177386106767, 12.46210287952087, 12.36330026954954, 12.32076967892521, 12.30501142428179, 12.30838980839569, 12.31946628099479, 12.3262961466662, 12.36612833040038, 12.39118256135069, 12.40582042984885, 12.41081561121008, 12.38849196828883, 12.3162034317334, 12.27098215056823, 12.29493249151987, 12.38970765691509, 12.50063075073655)
for ($i = 0; $i -lt $colH.Count; $i++) {
    $ws.Cells.Item(2 + $i, 8).Value = $colH[$i]
}

$colK = @(9.126514747940647, 8.810309577938671, 8.608831311384503, 8.524965713240354, 8.510936038122239, 8.607707289744983, 9.019053725272411, 9.764451962684868, 10.27115503628797, 10.49216350507358, 10.57444583537376, 10.55678815296652, 10.49896138544374, 10.46335615938052, 10.25651671757973, 10.12716301091572, 10.05187084093766, 10.02622654450662, 10.14102551885442, 10.51598506191226, 10.75281506901519, 10.62718022643571, 10.13476115686535, 9.569760534085882)
for ($i = 0; $i -lt $colK.Count; $i++) {
    $ws.Cells.Item(2 + $i, 11).Value = $colK[$i]
}

$colO = @(18.15487118542927, 18.27643367130139, 18.35587344356469, 18.38945098371154, 18.39509924875423, 18.35632140457804, 18.19578861013851, 17.91915182788767, 17.73929152551612, 17.66257389798908, 17.6342592299467, 17.64032448888341, 17.66022965223098, 17.67251816267789, 17.74440811667627, 17.78981904522569, 17.81641821964249, 17.82550659584685, 17.78493527380803, 17.65436300652959, 17.573321641638, 17.61618090300208, 17.78714169501942, 17.98989174822109)
for ($i = 0; $i -lt $colO.Count; $i++) {
    $ws.Cells.Item(2 + $i, 15).Value = $colO[$i]
}
